# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-34 (row 33 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 3
    9  = 0
    10 = 2
    11 = 2
    12 = 0
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 4
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 0
    24 = 0
    25 = 0
    26 = 2
    27 = 2
    28 = 2
    29 = 2
    30 = 1
    31 = 2
    32 = 1
    34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
